# Fixes issue of suzuki austrelia: fill in missing 0 values for
# section/header rows and a couple of partially-filled rows so that
# both the 2022 (C) and 2021 (D) columns are populated.

$wb = $excel.ActiveWorkbook

# --- Sheet "cbs_6" (Balance sheet) ---
$ws1 = $wb.Worksheets.Item("cbs_6")
$ws1.Range("C2").Value = 0
$ws1.Range("D2").Value = 0
$ws1.Range("C7").Value = 0
$ws1.Range("D7").Value = 0
$ws1.Range("C13").Value = 0
$ws1.Range("D13").Value = 0
$ws1.Range("C21").Value = 0
$ws1.Range("D21").Value = 0
$ws1.Range("C30").Value = 0
$ws1.Range("D30").Value = 0

# --- Sheet "cpl_5" (P&L) ---
$ws2 = $wb.Worksheets.Item("cpl_5")
$ws2.Range("C17").Value = 0
$ws2.Range("D17").Value = 0

# --- Sheet "ccf_8" (Cash flow) ---
$ws3 = $wb.Worksheets.Item("ccf_8")
$ws3.Range("C2").Value = 0
$ws3.Range("D2").Value = 0
$ws3.Range("C8").Value = 0
$ws3.Range("D8").Value = 0
$ws3.Range("C13").Value = 0
$ws3.Range("D13").Value = 0
$ws3.Range("C14").Value = 0
$ws3.Range("C15").Value = 0
